# Commit: "Added expirty to comments"
# Update timesheet values on the first worksheet ("Vārpas 1"):
#  - Row 2 (Nils Asejevs): D2 and E2 reset to 0, AE2 set to 2, AF2 set to 9,
#    with AG2 (Kopā) and AH2 (Dienas) recalculated accordingly.
#  - Row 3 (Signe Zalužinska): AF3 set to 9, with AG3 and AH3 recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 edits
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("AE2").Value = 2
$ws.Range("AF2").Value = 9
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 1.375

# Row 3 edits
$ws.Range("AF3").Value = 9
$ws.Range("AG3").Value = 9
$ws.Range("AH3").Value = 1.125
